# Updated symbol list on Wed Jan 25 14:32:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns D (Price) and E (Volume/1h) hold numeric/percent-looking text that
# must stay stored as text (matches the original inlineStr cells), so force
# the Text number format before writing the new value.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '300.76'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-4.23%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '35.44'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-1.61%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.048'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-1.16%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07980'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-1.98%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.911'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '-9.78%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '7.809'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '-1.75%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.9242'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-0.83%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1416'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '37.53%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1905'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-1.23%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.09212'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '1.76%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.03408'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-5.76%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.09877'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-0.14%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001387'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-3.11%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.005798'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '1.05%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.518'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '1.46%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.067'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-1.81%'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.983'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '3.37%'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3402'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-0.17%'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.1285'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-2.13%'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '-0.90%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.2405'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '8.49%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04503'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.001217'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-2.64%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.004785'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '13.43%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001233'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-1.55%'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0003006'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '-33.36%'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01909'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '-2.18%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04735'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '-3.12%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007342'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-3.48%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.009670'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '22.56%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1330'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-3.81%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.002115'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '0.39%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.01029'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-12.40%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006264'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-6.81%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000752'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '0.00%'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '57.72%'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.001662'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '-2.47%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002105'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.00%'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0002004'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '0.00%'
